# Daily attendance processing - 2025-11-24 11:47:50
# Reorders the "Recorded By" contributor lists for several sessions and
# refreshes the dependent attendance statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ANATOMY, Session 1
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# Row 3 - ANATOMY, Session 2
$ws.Range("G3").Value = "Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# Row 4 - ANATOMY, Session 3
$ws.Range("G4").Value = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg"

# Row 5 - ANATOMY, Session 4
$ws.Range("G5").Value = "eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"

# Row 7 - BIOCHEMISTRY LAB/CBL, Session 1
$ws.Range("G7").Value = "NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg"

# Row 9 - HISTOLOGY, Session 1
$ws.Range("G9").Value = "Safa.hany@med.asu.edu.eg, Shimaa.ashraf@med.asu.edu.eg"

# Row 10 - HISTOLOGY stats: Average Attendance %
# (written as literal text "26.0%", matching the source data which stores
# percentages as plain text rather than numeric percent values; a helper
# cell + copy/paste-values round trip avoids Excel's automatic text->number
# percent conversion while keeping the cell's existing style untouched)
$ws.Range("ZZ1").Formula = '="26.0%"'
$ws.Range("ZZ1").Copy()
$ws.Range("L10").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()

# Row 12 - MICROBIOLOGY, Session 1 (new recorders added, attendance count updated)
$ws.Range("G12").Value = "amira.m.ibrahim@med.asu.edu.eg, dina.adel@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H12").Value = "62/251"

# Row 15 - Group Statistics: Avg Attendance %
$ws.Range("ZZ1").Formula = '="26.0%"'
$ws.Range("ZZ1").Copy()
$ws.Range("S15").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
